$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Two newly solved questions (bit manipulation + tries) ---

# Row 9: "Min flips to make a OR b = c"
$ws.Range("A9").Value = 46073
$ws.Range("B9").Value = "Min flips to make a OR b = c"
$ws.Range("C9").Value = "https://leetcode.com/problems/minimum-flips-to-make-a-or-b-equal-to-c/description/"

# Row 10: "Implement Trie (Prefix Tree)"
$ws.Range("A10").Value = 46074
$ws.Range("B10").Value = "Implement Trie (Prefix Tree)"
$ws.Range("C10").Value = "https://leetcode.com/problems/implement-trie-prefix-tree/"

# Register the actual hyperlinks for the two new URL cells.
$ws.Hyperlinks.Add($ws.Range("C9"), "https://leetcode.com/problems/minimum-flips-to-make-a-or-b-equal-to-c/description/")
$ws.Hyperlinks.Add($ws.Range("C10"), "https://leetcode.com/problems/implement-trie-prefix-tree/")

# Copy the date formatting from the existing A8 cell down onto the two new
# date cells (A9, A10) so they pick up the same number formatting.
$ws.Range("A8").Copy()
$ws.Range("A9:A10").PasteSpecial(-4122)

# Normalize A8's own formatting to match the rest of the date column
# (A2:A7), which is what the author's commit actually did.
$ws.Range("A2").Copy()
$ws.Range("A8").PasteSpecial(-4122)

# Re-apply the same hyperlink-cell styling (underline/theme color) used by
# the rest of the URL column onto the two new cells -- Hyperlinks.Add()
# forces Excel's built-in "Hyperlink" style, so restore the sheet's own
# look here.
$ws.Range("C8").Copy()
$ws.Range("C9:C10").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# Drop the now-unused built-in "Hyperlink" named style that Hyperlinks.Add()
# registered, since nothing references it anymore.
$styles = $wb.Styles
for ($i = $styles.Count; $i -ge 1; $i--) {
    if ($styles.Item($i).Name -eq "Hyperlink") {
        $styles.Item($i).Delete()
    }
}

$wb.Save()
